$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("barnehage")

for ($row = 9; $row -le 35; $row++) {
    $ws.Cells.Item($row, 2).Value = $row - 1
}
